$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.029.70"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "3.226.56"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "3.213.08"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.482"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000235"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.20%  "
$ws.Range("D15").Value = "3.735.33"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "66.993.68"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "3.225.49"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("E18").Value = "  -2.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "502.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.718"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("E27").Value = "  -4.72%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "27.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.42%  "
$ws.Range("E31").Value = "  +4.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.72%  "
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "504.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0416"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0816"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.119"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("D42").Value = "2.853.87"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.252"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("E49").Value = "  -1.56%  "
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.19%  "
